$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save off the current "Terms Typically Offered" column (D) values
# before we overwrite column D with "Corequisites" data.
$termsHeader = $ws.Range("D1").Value2
$terms = @{}
for ($r = 2; $r -le 8; $r++) {
    $terms[$r] = $ws.Cells.Item($r, 4).Value2
}

# New header row: D=Corequisites, E=Concurrent, F=Recommended, G=Terms Typically Offered
$ws.Range("D1").Value2 = "Corequisites"
$ws.Range("E1").Value2 = "Concurrent"
$ws.Range("F1").Value2 = "Recommended"
$ws.Range("G1").Value2 = $termsHeader

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 4).Value2 = "NA"
    $ws.Cells.Item($r, 5).Value2 = "NA"
    $ws.Cells.Item($r, 6).Value2 = "NA"
    $ws.Cells.Item($r, 7).Value2 = $terms[$r]
}
